# Aula 18 - Planejamento Estrategico
# Fix: logo/background were only defined at the slideLayout/master level, which
# kept the deck from opening correctly. Embed the logo image + dark background
# individually in every slide and retire the now-unused slideLayout2 ("MASTER"
# layout) that used to carry them.

$p = $ppt.ActivePresentation

# Background color used on every slide (srgbClr 1A1A2E). COM RGB() packs as
# 0x00BBGGRR, so RGB(0x1A,0x1A,0x2E) == 0x2E1A1A.
$bgColor = 0x2E1A1A

# Logo geometry (EMU -> points, 1 pt = 12700 EMU):
#   off  x=274320  y=137160
#   ext  cx=1097280 cy=457200
$logoLeft   = 274320 / 12700
$logoTop    = 137160 / 12700
$logoWidth  = 1097280 / 12700
$logoHeight = 457200 / 12700

$defaultLayout = $p.SlideMaster.CustomLayouts.Item(1)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # 1) Dark navy background, set directly on the slide.
    $s.Background.Fill.Solid()
    $s.Background.Fill.ForeColor.RGB = $bgColor

    # 2) Logo image embedded on the slide itself (was previously only on the
    #    layout), sent behind every other shape so it sits first in z-order.
    $pic = $s.Shapes.AddPicture("preencoded.png", $false, $true, $logoLeft, $logoTop, $logoWidth, $logoHeight)
    $pic.Name = "Image 0"
    $pic.AlternativeText = "preencoded.png"
    $pic.LockAspectRatio = $true
    $pic.ZOrder(1)

    # 3) Move the slide off the soon-to-be-removed "MASTER" layout onto the
    #    plain default layout now that the slide owns its own bg + logo.
    $s.CustomLayout = $defaultLayout
}

# Drop the now-unused "MASTER" slide layout (index 2) that used to hold the
# shared background + logo; it is no longer referenced by any slide.
$unusedLayout = $p.SlideMaster.CustomLayouts.Item(2)
$unusedLayout.Delete()

# --- Slide-specific copy trims -------------------------------------------

$p.Slides.Item(2).Shapes.Item("Text 3").TextFrame.TextRange.Text = "Dominar 3 níveis"

$p.Slides.Item(3).Shapes.Item("Text 3").TextFrame.TextRange.Text = "Objetivo final"

$p.Slides.Item(4).Shapes.Item("Text 5").TextFrame.TextRange.Text = "Alocar blocos"
$p.Slides.Item(4).Shapes.Item("Text 7").TextFrame.TextRange.Text = "Buffers"

$p.Slides.Item(5).Shapes.Item("Text 3").TextFrame.TextRange.Text = "3 MITs"
$p.Slides.Item(5).Shapes.Item("Text 7").TextFrame.TextRange.Text = "Planejar noite anterior"

$p.Slides.Item(6).Shapes.Item("Text 1").TextFrame.TextRange.Text = "Revisão"
$p.Slides.Item(6).Shapes.Item("Text 5").TextFrame.TextRange.Text = "Ajustar plano"

$p.Slides.Item(7).Shapes.Item("Text 3").TextFrame.TextRange.Text = "3 níveis: longo, semanal, diário"
$p.Slides.Item(7).Shapes.Item("Text 5").TextFrame.TextRange.Text = "Reverse engineering"
$p.Slides.Item(7).Shapes.Item("Text 7").TextFrame.TextRange.Text = "MITs: 3 tarefas principais"
$p.Slides.Item(7).Shapes.Item("Text 9").TextFrame.TextRange.Text = "Revisar e ajustar"

$p.Slides.Item(8).Shapes.Item("Text 3").TextFrame.TextRange.Text = "Meu Mapa: objetivo + reverse + semanal + amanhã"

$p.Slides.Item(9).Shapes.Item("Text 3").TextFrame.TextRange.Text = '"Objetivo sem plano é apenas desejo."'
